# Append a new data row (row 41) with the latest Adafruit IO reading,
# mirroring the existing rows' layout: Timestamp | Feed Key | Value | Latitude | Longitude | Elevation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = "2024-09-25T18:06:40Z"
$ws.Range("B41").Value = "temperature"

# "25" must be stored as text (like the other Value cells), not a number.
# Temporarily force a text format so Excel doesn't auto-convert the numeric-looking
# string, then clear the format again so no stray cell style is left behind.
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "25"
$ws.Range("C41").ClearFormats()

$ws.Range("D41").Value = "N/A"
$ws.Range("E41").Value = "N/A"
$ws.Range("F41").Value = "N/A"
